# Update the S_SPT ant dimensions on row 5 (Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = 0.6
$ws.Range("C5").Value = 1.2
$ws.Range("D5").Value = 0.09
$ws.Range("E5").Value = 0.3
